# Automatische test-sync: 2025-07-29 21:33:50
$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new mail-log row (row 5) ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A5").Value = "Kun jij dit afhandelen?"
$ws.Range("B5").Value = "mailmind.test@zohomail.eu"
$ws.Range("C5").Value = "Testmail #3: Kun jij dit afhandelen?"
$ws.Range("D5").Value = "Overig"
$ws.Range("F5").Value = "2025-07-29 21:33:10"
$ws.Range("G5").Value = "Nee"
$ws.Range("H5").Value = "Ja"
$ws.Range("I5").Value = "Nee"
$ws.Range("J5").Value = "Nee"

# Extend the conditional-formatting ranges so they cover the new row too
$ws.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D5"))
$ws.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G5"))
$ws.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H5"))
$ws.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I5"))
$ws.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J5"))

# --- Dashboard sheet: bump the "Overig" tally from 2 to 3 ---
$ws2 = $wb.Worksheets.Item("Dashboard")
$ws2.Range("B2").Value = 3
